# Add new "Outliers_MAD" comparison columns (F, G, H) to the imputation
# comparison sheet, mirroring the existing KNN/SVM/RF columns (C, D, E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) --------------------------------------------------
# Copy the header formatting (bold, centered, bordered) from an existing
# header cell so the new headers match the look of ID/OriginalDataValue/etc.
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# --- Outlier flags (rows 2-18) -----------------------------------------
# Row => [KNN_Outliers_MAD, SVM_Outliers_MAD, RF_Outliers_MAD]
$outlierFlags = @{
    2  = @($false, $false, $true)
    3  = @($false, $false, $false)
    4  = @($false, $false, $false)
    5  = @($false, $false, $false)
    6  = @($false, $false, $false)
    7  = @($false, $false, $false)
    8  = @($false, $false, $false)
    9  = @($false, $false, $false)
    10 = @($false, $false, $false)
    11 = @($false, $false, $false)
    12 = @($false, $false, $false)
    13 = @($false, $false, $false)
    14 = @($false, $false, $false)
    15 = @($false, $false, $false)
    16 = @($false, $false, $true)
    17 = @($false, $false, $false)
    18 = @($false, $false, $false)
}

foreach ($row in $outlierFlags.Keys) {
    $vals = $outlierFlags[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
    $ws.Cells.Item($row, 8).Value = $vals[2]
}
